$p = $ppt.ActivePresentation

# Insert a new slide at position 2 (after the title slide, before the
# existing "MS, NMOSD, MOGAD" content slides), using the same
# Title-and-Content layout as the slide that is currently slide 2.
$layout = $p.Slides.Item(2).Layout
$s = $p.Slides.Add(2, $layout)

# Title placeholder
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Differential diagnosis of MS, NMOSD, MOGAD"

# Body / content placeholder - build up the bulleted outline
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Progressive demyelinating diseases of central nervous system (CNS)`rMultiple sclerosis`rNeuromyelitis optica spectrum disorder`rMyelin oligodendrocyte glycoprotein antibody associated disease`rOptic neuritis, myelitis, supra/infra-tentorial lesion`r"

$body.Paragraphs(2,1).IndentLevel = 2
$body.Paragraphs(3,1).IndentLevel = 2
$body.Paragraphs(4,1).IndentLevel = 2
